# Weekly data refresh for "Hortaliza, Vega Central Mapocho de Santiago - Acelga".
# A new week of price observations (date serial 44509 = 2021-11-09, qualities
# Extra/Primera/Segunda) is inserted at the top of the data block, pushing the
# existing rows 357:373 down to 360:376 (dimension grows from R373 to R376).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 357; rows 357:373 shift down to 360:376.
$ws.Rows.Item(357).Resize(3).Insert()

# New row 357 - Acelga, Extra
$ws.Range("A357").Value = 9
$ws.Range("B357").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C357").Value = "Metropolitana"
$ws.Range("D357").Value = 44509
$ws.Range("E357").Value = 13
$ws.Range("F357").Value = 100112009
$ws.Range("G357").Value = "Acelga"
$ws.Range("H357").Value = "Sin especificar"
$ws.Range("I357").Value = "Extra"
$ws.Range("J357").Value = 43
$ws.Range("K357").Value = 11000
$ws.Range("L357").Value = 11000
$ws.Range("M357").Value = 11000
$ws.Range("N357").Value = "$/docena de atados"
$ws.Range("O357").Value = "Región Metropolitana"
$ws.Range("P357").Value = 3667
$ws.Range("Q357").Value = 3
$ws.Range("R357").Value = "Hortaliza"

# New row 358 - Acelga, Primera
$ws.Range("A358").Value = 9
$ws.Range("B358").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C358").Value = "Metropolitana"
$ws.Range("D358").Value = 44509
$ws.Range("E358").Value = 13
$ws.Range("F358").Value = 100112009
$ws.Range("G358").Value = "Acelga"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 61
$ws.Range("K358").Value = 9000
$ws.Range("L358").Value = 10000
$ws.Range("M358").Value = 9508
$ws.Range("N358").Value = "$/docena de atados"
$ws.Range("O358").Value = "Región Metropolitana"
$ws.Range("P358").Value = 3169
$ws.Range("Q358").Value = 3
$ws.Range("R358").Value = "Hortaliza"

# New row 359 - Acelga, Segunda
$ws.Range("A359").Value = 9
$ws.Range("B359").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C359").Value = "Metropolitana"
$ws.Range("D359").Value = 44509
$ws.Range("E359").Value = 13
$ws.Range("F359").Value = 100112009
$ws.Range("G359").Value = "Acelga"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Segunda"
$ws.Range("J359").Value = 34
$ws.Range("K359").Value = 8000
$ws.Range("L359").Value = 8000
$ws.Range("M359").Value = 8000
$ws.Range("N359").Value = "$/docena de atados"
$ws.Range("O359").Value = "Región Metropolitana"
$ws.Range("P359").Value = 2667
$ws.Range("Q359").Value = 3
$ws.Range("R359").Value = "Hortaliza"
